# "Generate Report for Handback" - append a new handback row for
# 28017e47-eb36-4408-b9bc-7fea2d6061cd.md to the Overview / zh-cn / de-de
# sheets (mirrors the existing af758db3-... row already on each sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet -> new row 3
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(3, 1).Value = "28017e47-eb36-4408-b9bc-7fea2d6061cd.md"
$ov.Cells.Item(3, 3).Value = ".md"
$ov.Cells.Item(3, 5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3, 6).Value = "Handed back: in sync with en-US"

$ov.Cells.Item(3, 7).Value = "2016-10-14 07:38:19"
$ov.Cells.Item(3, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Cells.Item(3, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/241ab3f8dd3ac37ab8c934afd459b2653fca7e1b/e2e/28017e47-eb36-4408-b9bc-7fea2d6061cd.md", "", "", "e2e\28017e47-eb36-4408-b9bc-7fea2d6061cd.md")
$ov.Cells.Item(3, 2).Font.Underline = $true
$ov.Cells.Item(3, 2).Font.Color = 15570276

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G3"))

# ---------------------------------------------------------------------
# zh-cn sheet -> new row 3
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Cells.Item(3, 2).Value = ".md"
$zh.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(3, 4).Value = "e2e"
$zh.Cells.Item(3, 5).Value = "ht"
$zh.Cells.Item(3, 6).Value = "True"
$zh.Cells.Item(3, 7).Value = "28017e47-eb36-4408-b9bc-7fea2d6061cd.9b2e2be9b89b3d9836bf682be1bfe38e9a0d79f4.zh-cn.xlf"
$zh.Cells.Item(3, 10).Value = "28017e47-eb36-4408-b9bc-7fea2d6061cd.9b2e2be9b89b3d9836bf682be1bfe38e9a0d79f4.zh-cn.xlf"
$zh.Cells.Item(3, 13).Value = "True"
$zh.Cells.Item(3, 15).Value = "False"

$zh.Cells.Item(3, 8).Value = "2016-10-14 07:38:09"
$zh.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item(3, 11).Value = "2016-10-14 07:38:52"
$zh.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zh.Hyperlinks.Add($zh.Cells.Item(3, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/241ab3f8dd3ac37ab8c934afd459b2653fca7e1b/e2e/28017e47-eb36-4408-b9bc-7fea2d6061cd.md", "", "", "28017e47-eb36-4408-b9bc-7fea2d6061cd.md")
$zh.Cells.Item(3, 1).Font.Underline = $true
$zh.Cells.Item(3, 1).Font.Color = 15570276

$zh.Hyperlinks.Add($zh.Cells.Item(3, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7f6134260f428ae26ee75c7e32c2bab346fc3ef3/e2e/28017e47-eb36-4408-b9bc-7fea2d6061cd.md", "", "", "28017e47-eb36-4408-b9bc-7fea2d6061cd.md")
$zh.Cells.Item(3, 9).Font.Underline = $true
$zh.Cells.Item(3, 9).Font.Color = 15570276

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P3"))

# ---------------------------------------------------------------------
# de-de sheet -> new row 3
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Cells.Item(3, 2).Value = ".md"
$de.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
$de.Cells.Item(3, 4).Value = "e2e"
$de.Cells.Item(3, 5).Value = "ht"
$de.Cells.Item(3, 6).Value = "True"
$de.Cells.Item(3, 7).Value = "28017e47-eb36-4408-b9bc-7fea2d6061cd.9b2e2be9b89b3d9836bf682be1bfe38e9a0d79f4.de-de.xlf"
$de.Cells.Item(3, 10).Value = "28017e47-eb36-4408-b9bc-7fea2d6061cd.9b2e2be9b89b3d9836bf682be1bfe38e9a0d79f4.de-de.xlf"
$de.Cells.Item(3, 13).Value = "True"
$de.Cells.Item(3, 15).Value = "False"

$de.Cells.Item(3, 8).Value = "2016-10-14 07:38:19"
$de.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item(3, 11).Value = "2016-10-14 07:39:09"
$de.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$de.Hyperlinks.Add($de.Cells.Item(3, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/241ab3f8dd3ac37ab8c934afd459b2653fca7e1b/e2e/28017e47-eb36-4408-b9bc-7fea2d6061cd.md", "", "", "28017e47-eb36-4408-b9bc-7fea2d6061cd.md")
$de.Cells.Item(3, 1).Font.Underline = $true
$de.Cells.Item(3, 1).Font.Color = 15570276

$de.Hyperlinks.Add($de.Cells.Item(3, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/76bfa6688c4e7ed2b2a16fee5790e3418e16e474/e2e/28017e47-eb36-4408-b9bc-7fea2d6061cd.md", "", "", "28017e47-eb36-4408-b9bc-7fea2d6061cd.md")
$de.Cells.Item(3, 9).Font.Underline = $true
$de.Cells.Item(3, 9).Font.Color = 15570276

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P3"))
